$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 522-523, pushing the existing rows 522..581
# down to 524..583 (dimension grows from A1:R581 to A1:R583).
$ws.Range("A522:A523").EntireRow.Insert()

# New row 522 - Zafiro rojo entry dated 45132
$ws.Cells.Item(522,1).Value = 11
$ws.Cells.Item(522,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(522,3).Value = "Bíobío"
$ws.Cells.Item(522,4).Value = 45132
$ws.Cells.Item(522,5).Value = 8
$ws.Cells.Item(522,6).Value = 100112002
$ws.Cells.Item(522,7).Value = "Pimiento"
$ws.Cells.Item(522,8).Value = "Zafiro rojo"
$ws.Cells.Item(522,9).Value = "Primera"
$ws.Cells.Item(522,10).Value = 100
$ws.Cells.Item(522,11).Value = 19000
$ws.Cells.Item(522,12).Value = 20000
$ws.Cells.Item(522,13).Value = 19500
$ws.Cells.Item(522,14).Value = "$/caja 15 kilos"
$ws.Cells.Item(522,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(522,16).Value = 1300
$ws.Cells.Item(522,17).Value = 15
$ws.Cells.Item(522,18).Value = "Hortaliza"

# New row 523 - Zafiro verde entry dated 45132
$ws.Cells.Item(523,1).Value = 11
$ws.Cells.Item(523,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(523,3).Value = "Bíobío"
$ws.Cells.Item(523,4).Value = 45132
$ws.Cells.Item(523,5).Value = 8
$ws.Cells.Item(523,6).Value = 100112002
$ws.Cells.Item(523,7).Value = "Pimiento"
$ws.Cells.Item(523,8).Value = "Zafiro verde"
$ws.Cells.Item(523,9).Value = "Primera"
$ws.Cells.Item(523,10).Value = 100
$ws.Cells.Item(523,11).Value = 15000
$ws.Cells.Item(523,12).Value = 16000
$ws.Cells.Item(523,13).Value = 15500
$ws.Cells.Item(523,14).Value = "$/caja 15 kilos"
$ws.Cells.Item(523,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(523,16).Value = 1033
$ws.Cells.Item(523,17).Value = 15
$ws.Cells.Item(523,18).Value = "Hortaliza"

Write-Host "Inserted 2 rows; UsedRange rows:" $ws.UsedRange.Rows.Count
